# Commit: "Added 02 call vs put smile plot"
# Inserts a new "Dividends" worksheet (CAT dividend schedule pulled from
# Bloomberg, plus a short dividend-yield-assumption note) between the
# existing "ATM" and "ATM Option" tabs.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert the new worksheet immediately before "ATM Option" so the
#    tab order becomes: Calls, Puts, Price, USGG1M, ATM, Dividends,
#    ATM Option, MetaData.
# ------------------------------------------------------------------
$atmOption = $wb.Worksheets.Item("ATM Option")
$ws = $wb.Worksheets.Add($atmOption)
$ws.Name = "Dividends"

# A source cell elsewhere in the workbook that already carries the
# plain date number format (style reused instead of minting a new one).
$dateFmtSource = $wb.Worksheets.Item("Price").Range("A2")

# ------------------------------------------------------------------
# 2. Header row.
# ------------------------------------------------------------------
$ws.Range("A1").Value = "Declaration"
$ws.Range("B1").Value = "Ex Date"
$ws.Range("C1").Value = "Record"
$ws.Range("D1").Value = "Payable"
$ws.Range("E1").Value = "Curr"
$ws.Range("F1").Value = "Amount"
$ws.Range("G1").Value = "Type"

# ------------------------------------------------------------------
# 3. Dividend schedule (Bloomberg BDVD export) - apply date format to
#    the four date columns first, then fill in the values.
# ------------------------------------------------------------------
$dateFmtSource.Copy()
$ws.Range("A2:D10").PasteSpecial(-4122)

$rows = @(
    @("12/10/2025","01/20/2026","10/20/2025","11/20/2025","USD",1.51,"BDVD Projection"),
    @("10/06/2025","10/20/2025","10/20/2025","11/20/2025","USD",1.51,"Regular Cash"),
    @("06/11/2025","07/21/2025","07/21/2025","08/20/2025","USD",1.51,"Regular Cash"),
    @("04/09/2025","04/21/2025","04/21/2025","05/20/2025","USD",1.41,"Regular Cash"),
    @("12/10/2024","01/21/2025","01/21/2025","02/20/2025","USD",1.41,"Regular Cash"),
    @("10/09/2024","10/21/2024","10/21/2024","11/20/2024","USD",1.41,"Regular Cash"),
    @("06/12/2024","07/22/2024","07/22/2024","08/20/2024","USD",1.41,"Regular Cash"),
    @("04/10/2024","04/19/2024","04/22/2024","05/20/2024","USD",1.30,"Regular Cash"),
    @("12/13/2023","01/19/2024","01/22/2024","02/20/2024","USD",1.30,"Regular Cash")
)

$r = 2
foreach ($row in $rows) {
    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
    $ws.Range("E$r").Value = $row[4]
    $ws.Range("F$r").Value = $row[5]
    $ws.Range("G$r").Value = $row[6]
    $r++
}

# ------------------------------------------------------------------
# 4. Dividend-yield-assumption note below the table. Rows 12-13 and 15
#    keep the (empty) date-formatted style left over from extending the
#    format down column A; rows 16-20 carry the note text in that same
#    style, and the final block (21-26) is plain/unformatted.
# ------------------------------------------------------------------
$dateFmtSource.Copy()
$ws.Range("A12:A20").PasteSpecial(-4122)

$ws.Range("H21").Value = " "
$ws.Range("A16").Value = "We obtained the CAT dividend schedule from Bloomberg. "
$ws.Range("A17").Value = "The ex dates around our trade date (19 Sep 2025) are:"
$ws.Range("A19").Value = "- 21 Jul 2025 (already passed), and "
$ws.Range("A20").Value = "- 20 Oct 2025 (after the option expiry 17 Oct 2025)."
$ws.Range("A22").Value = "Hence there are no cash dividends paid between the trade date and the "
$ws.Range("A23").Value = "option maturity. For the life of this option the stock behaves as a "
$ws.Range("A24").Value = "non dividend paying asset, so the standard Black Scholes and put–call "
$ws.Range("A25").Value = "parity formulas without dividends are appropriate. We therefore set "
$ws.Range("A26").Value = "the dividend yield `$q = 0`$."
$ws.Range("A14").Value = "Dividend assumption"

# ------------------------------------------------------------------
# 5. Column widths, matching the widths Excel settled on after the
#    author auto-fit the date/label columns to their content (the long
#    note in column A keeps the narrower, pre-existing width rather
#    than stretching to fit that wrapped paragraph).
# ------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 9.8
$ws.Columns.Item(2).ColumnWidth = 9.62
$ws.Columns.Item(3).ColumnWidth = 9.62
$ws.Columns.Item(4).ColumnWidth = 9.62
$ws.Columns.Item(5).ColumnWidth = 3.62
$ws.Columns.Item(6).ColumnWidth = 6.89
$ws.Columns.Item(7).ColumnWidth = 14.17
$ws.Columns.Item(8).ColumnWidth = 0.53

# ------------------------------------------------------------------
# 6. Leave the cursor where the author left it and make this the
#    active tab (mirrors tabSelected moving off "ATM").
# ------------------------------------------------------------------
$ws.Range("A15").Select()
$ws.Activate()
